$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value = 296.9
$ws.Range("G7").Value = 303.15
$ws.Range("H7").Value = 295.3
$ws.Range("I7").Value = 301.95
$ws.Range("J7").Value = 297.65

$ws.Range("G9").Value = 299
$ws.Range("H9").Value = 293.2
$ws.Range("I9").Value = 297.2

$ws.Range("G10").Value = 299.5
$ws.Range("H10").Value = 296.15
$ws.Range("I10").Value = 296.75

$ws.Range("G11").Value = 300.6
$ws.Range("H11").Value = 296.3
$ws.Range("I11").Value = 299.65

$ws.Range("G12").Value = 299.75
$ws.Range("H12").Value = 296.85
$ws.Range("I12").Value = 298

$ws.Range("G13").Value = 299.2
$ws.Range("H13").Value = 297
$ws.Range("I13").Value = 297.35

$ws.Range("G14").Value = 298.65
$ws.Range("H14").Value = 297.1
$ws.Range("I14").Value = 298.25

$ws.Range("G15").Value = 299.65
$ws.Range("H15").Value = 297.9
$ws.Range("I15").Value = 298.65

$ws.Range("G16").Value = 302.1
$ws.Range("H16").Value = 298.5
$ws.Range("I16").Value = 300.75

$ws.Range("G17").Value = 301.5
$ws.Range("H17").Value = 299.65
$ws.Range("I17").Value = 301

$ws.Range("G18").Value = 301.45
$ws.Range("H18").Value = 298.6
$ws.Range("I18").Value = 299.75

$ws.Range("G19").Value = 300.75
$ws.Range("H19").Value = 298
$ws.Range("I19").Value = 299.95

$ws.Range("G20").Value = 301.7
$ws.Range("H20").Value = 299.9
$ws.Range("I20").Value = 301.6

$ws.Range("G21").Value = 303.15
$ws.Range("H21").Value = 300.9
$ws.Range("I21").Value = 302.45
